$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B3: "3-item survey scale (COVIDiStress); ... First Stage Dependent Variable" ---
$b3text1 = "3-item survey scale (COVIDiStress);                              "
$b3text2 = "First Stage Dependent Variable"
$ws.Range("B3").Value = $b3text1 + $b3text2
$b3run = $ws.Range("B3").Characters($b3text1.Length + 1, $b3text2.Length)
$b3run.Font.Italic = $true

# --- B8: "Ratio of infection rate, May 31st to 1st, ... Second Stage Dependent Variable" ---
$b8text1 = "Ratio of infection rate, May 31st to 1st, (Johns Hopkins, 18-day lead in COVID-19 deaths); "
$b8text2 = "Second Stage Dependent Variable"
$ws.Range("B8").Value = $b8text1 + $b8text2
$b8run = $ws.Range("B8").Characters($b8text1.Length + 1, $b8text2.Length)
$b8run.Font.Italic = $true

# --- Bold the two variable names that are now called out as the IV / DV ---
$ws.Range("A3").Font.Bold = $true
$ws.Range("A8").Font.Bold = $true
